$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at 88-89, shifting existing rows 88-104 down to 90-106.
$ws.Rows("88:89").Insert()

# Row 88 - new weekly entry (Primera)
$ws.Range("A88").Value = 1
$ws.Range("B88").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C88").Value = "Arica y Parinacota"
$ws.Range("D88").Value = 44722
$ws.Range("E88").Value = 15
$ws.Range("F88").Value = 100112036
$ws.Range("G88").Value = "Caigua"
$ws.Range("H88").Value = "Sin especificar"
$ws.Range("I88").Value = "Primera"
$ws.Range("J88").Value = 120
$ws.Range("K88").Value = 9000
$ws.Range("L88").Value = 10000
$ws.Range("M88").Value = 9500
$ws.Range("N88").Value = "$/caja 20 kilos"
$ws.Range("O88").Value = "Región de Arica y Parinacota"
$ws.Range("P88").Value = 475
$ws.Range("Q88").Value = 20
$ws.Range("R88").Value = "Hortaliza"

# Row 89 - new weekly entry (Segunda)
$ws.Range("A89").Value = 1
$ws.Range("B89").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C89").Value = "Arica y Parinacota"
$ws.Range("D89").Value = 44722
$ws.Range("E89").Value = 15
$ws.Range("F89").Value = 100112036
$ws.Range("G89").Value = "Caigua"
$ws.Range("H89").Value = "Sin especificar"
$ws.Range("I89").Value = "Segunda"
$ws.Range("J89").Value = 160
$ws.Range("K89").Value = 7000
$ws.Range("L89").Value = 8000
$ws.Range("M89").Value = 7500
$ws.Range("N89").Value = "$/caja 20 kilos"
$ws.Range("O89").Value = "Región de Arica y Parinacota"
$ws.Range("P89").Value = 375
$ws.Range("Q89").Value = 20
$ws.Range("R89").Value = "Hortaliza"
